$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update sheet (tab) name to reflect new "through" date
$ws.Name = "Through 2021-10-10"

# Update header label for October row
$ws.Range("A11").Value = "October (through 10-10)"

# Update October row (row 11) values
$ws.Range("B11").Value = 8
$ws.Range("C11").Value = 17
$ws.Range("D11").Value = 17
$ws.Range("E11").Value = 25
$ws.Range("F11").Value = 11
$ws.Range("G11").Value = 43
$ws.Range("H11").Value = 66

# Update Total row (row 12) values
$ws.Range("B12").Value = 234
$ws.Range("C12").Value = 446
$ws.Range("D12").Value = 644
$ws.Range("E12").Value = 573
$ws.Range("F12").Value = 433
$ws.Range("G12").Value = 944
$ws.Range("H12").Value = 1316
